$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for columns D, L, M, N, O, P, Q, R, S, T (rows 2-20)
# since the edit is a permutation of these values across rows (the other
# columns A,B,C,E,F,G,H,I,J,K are identical on every row, so they need no change).
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: new row -> source (old) row whose values it should receive
$perm = @{
    2 = 13
    3 = 19
    4 = 17
    5 = 3
    6 = 20
    7 = 2
    8 = 14
    9 = 11
    10 = 16
    11 = 12
    12 = 8
    13 = 6
    14 = 18
    15 = 15
    16 = 7
    17 = 4
    18 = 5
    19 = 10
    20 = 9
}

foreach ($destRow in ($perm.Keys | Sort-Object)) {
    $srcRow = $perm[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
